# Updated RAD Test Data for Existing Liability MD CRN.
#
# The "Existing" worksheet gains a new "CRN" column between the
# FeinSsn/FEINSSN column and the "App ID" column. The old "App ID" /
# "Backend Data" columns shift one column to the right, and several rows
# get a new "Y" flag in the inserted CRN column. The "Existing" sheet
# also becomes the active sheet/tab (previously "Personal_EL" was active).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Existing")

# Insert a new blank column at G (before the old "App ID" column), which
# shifts "App ID" -> H and "Backend Data" -> I, extending the used range
# from A1:H19 to A1:I19.
$ws.Columns.Item(7).Insert() | Out-Null

# Inserting a column copies the formatting of the column to its left (F)
# into any row where F had an explicit style, leaving a few empty-but-
# styled cells behind in the new column. Clear those so the rows that
# should have no CRN entry come back out completely blank (no <c> at
# all), matching the original author's edit.
$emptyCrnRows = @(5, 6, 11, 12)
foreach ($r in $emptyCrnRows) {
    $ws.Cells.Item($r, 7).Clear() | Out-Null
}

# New header for the inserted column.
$ws.Range("G1").Value = "CRN"

# Rows that get a "Y" flag in the new CRN column.
$crnYesRows = @(7, 8, 9, 10, 13, 14, 16, 17)
foreach ($r in $crnYesRows) {
    $ws.Cells.Item($r, 7).Value = "Y"
}

# Make "Existing" the active sheet/tab and set the active selection to
# D23 (a cell below the data, matching the author's saved view).
$ws.Activate()
$ws.Range("D23").Select() | Out-Null
